$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.101478576660156
$ws.Range("B1").Value = 0.8685930967330933
$ws.Range("C1").Value = 3.265453338623047
$ws.Range("D1").Value = 3.165665864944458
$ws.Range("E1").Value = 0.9391952753067017
